# "added queue to filereader": reshuffles the course-schedule class lists
# across the 08:00-14:15 time slots (rows 2-7, columns B-F).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

# 08:00 row
$ws.Range("B2").Value = "CMPS210${nl}ENG205${nl}PHYS110${nl}"
$ws.Range("C2").Value = "MATH201${nl}HIST301${nl}ART202${nl}BUSN301${nl}"
$ws.Range("D2").Value = "CMPS210${nl}ENG205${nl}PHYS110${nl}"
$ws.Range("E2").Value = "MATH201${nl}HIST301${nl}ART202${nl}BUSN301${nl}"
$ws.Range("F2").Value = "PSYC301${nl}PHYS110${nl}LANG202${nl}"

# 09:15 row
$ws.Range("C3").Value = "MATH201${nl}HIST301${nl}ART202${nl}BUSN301${nl}"
$ws.Range("D3").Value = "ENG205${nl}PHYS110${nl}"
$ws.Range("E3").Value = "MATH201${nl}HIST301${nl}ART202${nl}BUSN301${nl}"
$ws.Range("F3").Value = "PHYS110${nl}LANG202${nl}"

# 10:30 row
$ws.Range("C4").Value = "MATH201${nl}"
$ws.Range("E4").Value = "MATH201${nl}"

# 11:45 row
$ws.Range("B5").Value = "PSYC301${nl}"
$ws.Range("C5").Value = "CHEM202${nl}"
$ws.Range("D5").Value = "PSYC301${nl}"
$ws.Range("E5").Value = "CHEM202${nl}"
$ws.Range("F5").Value = "PSYC301${nl}"

# 13:00 row
$ws.Range("D6").Value = "PSYC301${nl}"

# 14:15 row
$ws.Range("C7").Value = "CHEM202${nl}"
$ws.Range("E7").Value = "CHEM202${nl}"
